$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'tropfen',
    'bitten',
    'werden',
    'sprengen',
    'hören',
    'räumen',
    'lesen',
    'scheinen',
    'freuen',
    'sinken',
    'wenden',
    'feiern',
    'greifen',
    'runden',
    'töten',
    'stehlen',
    'werfen',
    'erben',
    'heben',
    'dienen',
    'suchen',
    'reizen',
    'schmecken',
    'brauchen',
    'malen',
    'geben',
    'trauen',
    'tollen',
    'ehren',
    'zeigen',
    'fahren',
    'liegen',
    'zögern',
    'schenken',
    'schlucken',
    'fällen',
    'wehtun',
    'spüren',
    'platzen',
    'quälen',
    'zielen',
    'hauen',
    'biegen',
    'grüßen',
    'irren',
    'stammen',
    'arten',
    'backen',
    'filmen',
    'pflanzen',
    'fangen',
    'kehren',
    'schreiten',
    'dringen',
    'flüchten',
    'gelten',
    'heulen',
    'wüten',
    'bergen',
    'seufzen',
    'rasen',
    'spinnen',
    'drehen',
    'formen',
    'streichen',
    'klingen',
    'betteln',
    'lügen',
    'knarren',
    'schrecken',
    'warnen',
    'sorgen',
    'graben',
    'mögen',
    'jubeln',
    'kichern',
    'sperren',
    'schwören',
    'stecken',
    'kosten',
    'pfeifen',
    'spielen',
    'münzen',
    'heilen',
    'siegen',
    'helfen',
    'jagen',
    'weichen',
    'bauen',
    'gründen',
    'trennen',
    'wachsen',
    'saufen',
    'ärgern',
    'folgen',
    'treiben',
    'kümmern',
    'schulden',
    'achten',
    'bellen',
    'machen',
    'flehen',
    'liefern',
    'mauern',
    'führen',
    'scheitern',
    'fallen',
    'planen',
    'wirken',
    'lockern',
    'fließen',
    'schwingen',
    'bluten',
    'zünden',
    'ändern',
    'fischen',
    'enden',
    'äußern',
    'sichern',
    'wundern',
    'klettern',
    'altern',
    'boxen',
    'loben',
    'decken',
    'rufen',
    'sterben',
    'schwächen'
)

for ($i = 0; $i -lt $words.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $words[$i]
}

